# Auto-generated edit script: updates cached market-price / profit
# values across the Alpha_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, WVR)
# to match a refreshed data pull from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1480.7693
$ws.Range("I9").Value = 95
$ws.Range("J9").Value = 2346.875
$ws.Range("K9").Value = 95
$ws.Range("L9").Value = 2346.875
$ws.Range("M9").Value = 74
$ws.Range("N9").Value = -2684.875
$ws.Range("H12").Value = 358.75
$ws.Range("I12").Value = 338.33334
$ws.Range("J12").Value = 420
$ws.Range("K12").Value = 338.33334
$ws.Range("L12").Value = 420
$ws.Range("M12").Value = -168.33334
$ws.Range("N12").Value = -760
$ws.Range("H19").Value = 629
$ws.Range("I19").Value = 500
$ws.Range("K19").Value = 500
$ws.Range("M19").Value = -325
$ws.Range("H98").Value = 2621.75
$ws.Range("I98").Value = 2843.9565
$ws.Range("J98").Value = 1599.6
$ws.Range("K98").Value = 2843.9565
$ws.Range("L98").Value = 1599.6
$ws.Range("M98").Value = -1345.9565
$ws.Range("N98").Value = -4595.6
$ws.Range("H107").Value = 2582
$ws.Range("I107").Value = 2873
$ws.Range("K107").Value = 2873
$ws.Range("M107").Value = -953
$ws.Range("H122").Value = 2621.75
$ws.Range("I122").Value = 2843.9565
$ws.Range("J122").Value = 1599.6
$ws.Range("K122").Value = 8531.869499999999
$ws.Range("L122").Value = 4798.799999999999
$ws.Range("M122").Value = -6081.869499999999
$ws.Range("N122").Value = -9698.799999999999
$ws.Range("H129").Value = 2008.6263
$ws.Range("J129").Value = 2097.2092
$ws.Range("L129").Value = 6291.6276
$ws.Range("N129").Value = -16291.6276
$ws.Range("H132").Value = 1499.8684
$ws.Range("I132").Value = 1374.8889
$ws.Range("J132").Value = 3749.5
$ws.Range("K132").Value = 4124.6667
$ws.Range("L132").Value = 11248.5
$ws.Range("M132").Value = -1594.6667
$ws.Range("N132").Value = -16308.5
$ws.Range("H137").Value = 2287.2778
$ws.Range("I137").Value = 1759
$ws.Range("K137").Value = 5277
$ws.Range("M137").Value = -2727
$ws.Range("H138").Value = 3151.8086
$ws.Range("J138").Value = 3824.2222
$ws.Range("L138").Value = 11472.6666
$ws.Range("N138").Value = -21752.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2764.6667
$ws.Range("I74").Value = 2440
$ws.Range("K74").Value = 2440
$ws.Range("M74").Value = -1566
$ws.Range("H77").Value = 2764.6667
$ws.Range("I77").Value = 2440
$ws.Range("K77").Value = 12200
$ws.Range("M77").Value = -7832
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H98").Value = 13666
$ws.Range("J98").Value = 13666
$ws.Range("L98").Value = 13666
$ws.Range("N98").Value = -19656
$ws.Range("H101").Value = 19949.5
$ws.Range("J101").Value = 19949.5
$ws.Range("L101").Value = 19949.5
$ws.Range("N101").Value = -26439.5
$ws.Range("H102").Value = 1431.8
$ws.Range("I102").Value = 1037
$ws.Range("K102").Value = 1037
$ws.Range("M102").Value = 585
$ws.Range("H112").Value = 12747.167
$ws.Range("J112").Value = 12747.167
$ws.Range("L112").Value = 12747.167
$ws.Range("N112").Value = -15701.167

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 359
$ws.Range("I22").Value = 348.75
$ws.Range("K22").Value = 348.75
$ws.Range("M22").Value = -175.75
$ws.Range("H94").Value = 1550.1765
$ws.Range("I94").Value = 1428.7693
$ws.Range("J94").Value = 1944.75
$ws.Range("K94").Value = 1428.7693
$ws.Range("L94").Value = 1944.75
$ws.Range("M94").Value = -977.7692999999999
$ws.Range("N94").Value = -2846.75
$ws.Range("H105").Value = 976
$ws.Range("I105").Value = 981.6
$ws.Range("J105").Value = 966.6667
$ws.Range("K105").Value = 981.6
$ws.Range("L105").Value = 966.6667
$ws.Range("M105").Value = 765.4
$ws.Range("N105").Value = -4460.6667
$ws.Range("H107").Value = 3571.7144
$ws.Range("I107").Value = 2200.4
$ws.Range("J107").Value = 7000
$ws.Range("K107").Value = 2200.4
$ws.Range("L107").Value = 7000
$ws.Range("M107").Value = -280.4000000000001
$ws.Range("N107").Value = -10840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 39000
$ws.Range("J28").Value = 39000
$ws.Range("L28").Value = 39000
$ws.Range("N28").Value = -39490
$ws.Range("H31").Value = 2794.9285
$ws.Range("I31").Value = 2745
$ws.Range("K31").Value = 2745
$ws.Range("M31").Value = -2450
$ws.Range("H34").Value = 2794.9285
$ws.Range("I34").Value = 2745
$ws.Range("K34").Value = 2745
$ws.Range("M34").Value = -2543
$ws.Range("H58").Value = 1829.3334
$ws.Range("I58").Value = 1332.3334
$ws.Range("K58").Value = 1332.3334
$ws.Range("M58").Value = -1129.3334
$ws.Range("H62").Value = 41056.062
$ws.Range("I62").Value = 3876.6924
$ws.Range("J62").Value = 202166.67
$ws.Range("K62").Value = 3876.6924
$ws.Range("L62").Value = 202166.67
$ws.Range("M62").Value = -3252.6924
$ws.Range("N62").Value = -203414.67
$ws.Range("H65").Value = 41056.062
$ws.Range("I65").Value = 3876.6924
$ws.Range("J65").Value = 202166.67
$ws.Range("K65").Value = 19383.462
$ws.Range("L65").Value = 1010833.35
$ws.Range("M65").Value = -16263.462
$ws.Range("N65").Value = -1017073.35
$ws.Range("H96").Value = 14970.143
$ws.Range("J96").Value = 14970.143
$ws.Range("L96").Value = 14970.143
$ws.Range("N96").Value = -20462.143
$ws.Range("H99").Value = 4899.6
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 4899.6
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 4899.6
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -7895.6
$ws.Range("H105").Value = 1277.4445
$ws.Range("I105").Value = 879.5
$ws.Range("J105").Value = 2073.3333
$ws.Range("K105").Value = 879.5
$ws.Range("L105").Value = 2073.3333
$ws.Range("M105").Value = 867.5
$ws.Range("N105").Value = -5567.3333
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H126").Value = 4899.6
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4899.6
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 14698.8
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -19638.8
$ws.Range("H136").Value = 1829.3334
$ws.Range("I136").Value = 1332.3334
$ws.Range("K136").Value = 3997.0002
$ws.Range("M136").Value = -1447.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 2897.6
$ws.Range("I13").Value = 2996
$ws.Range("K13").Value = 8988
$ws.Range("M13").Value = -8820
$ws.Range("H23").Value = 175.6
$ws.Range("I23").Value = 173.4
$ws.Range("J23").Value = 177.8
$ws.Range("K23").Value = 520.2
$ws.Range("L23").Value = 533.4000000000001
$ws.Range("M23").Value = -285.2
$ws.Range("N23").Value = -1003.4
$ws.Range("H132").Value = 5384.75
$ws.Range("I132").Value = 1766.2222
$ws.Range("J132").Value = 10037.143
$ws.Range("K132").Value = 15895.9998
$ws.Range("L132").Value = 90334.287
$ws.Range("M132").Value = -13365.9998
$ws.Range("N132").Value = -95394.287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 20000000
$ws.Range("I10").Value = 20000000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 20000000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -19999831
$ws.Range("N10").ClearContents()
$ws.Range("H80").Value = 4905.643
$ws.Range("I80").Value = 2959.3333
$ws.Range("J80").Value = 5436.4546
$ws.Range("K80").Value = 2959.3333
$ws.Range("L80").Value = 5436.4546
$ws.Range("M80").Value = -1961.3333
$ws.Range("N80").Value = -7432.4546
$ws.Range("H83").Value = 4905.643
$ws.Range("I83").Value = 2959.3333
$ws.Range("J83").Value = 5436.4546
$ws.Range("K83").Value = 14796.6665
$ws.Range("L83").Value = 27182.273
$ws.Range("M83").Value = -9804.666499999999
$ws.Range("N83").Value = -37166.273
$ws.Range("H102").Value = 1204.8334
$ws.Range("J102").Value = 1289
$ws.Range("L102").Value = 1289
$ws.Range("N102").Value = -4533
$ws.Range("H107").Value = 531
$ws.Range("I107").Value = 247.88889
$ws.Range("J107").Value = 1380.3334
$ws.Range("K107").Value = 247.88889
$ws.Range("L107").Value = 1380.3334
$ws.Range("M107").Value = 1672.11111
$ws.Range("N107").Value = -5220.3334
$ws.Range("H132").Value = 2290.9
$ws.Range("I132").Value = 1970.5
$ws.Range("J132").Value = 2504.5
$ws.Range("K132").Value = 5911.5
$ws.Range("L132").Value = 7513.5
$ws.Range("M132").Value = -3381.5
$ws.Range("N132").Value = -12573.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 28635.781
$ws.Range("I15").Value = 9999
$ws.Range("K15").Value = 9999
$ws.Range("M15").Value = -9711
$ws.Range("H34").Value = 24996
$ws.Range("I34").Value = 27999.5
$ws.Range("J34").Value = 21992.5
$ws.Range("K34").Value = 27999.5
$ws.Range("L34").Value = 21992.5
$ws.Range("M34").Value = -27796.5
$ws.Range("N34").Value = -22398.5
$ws.Range("H74").Value = 15951.25
$ws.Range("J74").Value = 18468.666
$ws.Range("L74").Value = 18468.666
$ws.Range("N74").Value = -20340.666
$ws.Range("H77").Value = 15951.25
$ws.Range("J77").Value = 18468.666
$ws.Range("L77").Value = 55405.99800000001
$ws.Range("N77").Value = -64765.99800000001
$ws.Range("H126").Value = 2045.6842
$ws.Range("I126").Value = 2025.3334
$ws.Range("J126").Value = 2122
$ws.Range("K126").Value = 6076.0002
$ws.Range("L126").Value = 6366
$ws.Range("M126").Value = -3606.0002
$ws.Range("N126").Value = -11306
